# Insert a new weekly record row at row 141, pushing existing rows 141:240
# down to 142:241 (dimension grows from A1:R240 to A1:R241).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(141).Insert()

$ws.Range("A141").Value = 4
$ws.Range("B141").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value = "Los Lagos"
$ws.Range("D141").Value = Get-Date -Year 2022 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Range("E141").Value = 10
$ws.Range("F141").Value = 100112021
$ws.Range("G141").Value = "Ají"
$ws.Range("H141").Value = "Inferno"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 40
$ws.Range("K141").Value = 22000
$ws.Range("L141").Value = 22000
$ws.Range("M141").Value = 22000
$ws.Range("N141").Value = '$/caja 12 kilos'
$ws.Range("O141").Value = "Región de Arica y Parinacota"
$ws.Range("P141").Value = 1833
$ws.Range("Q141").Value = 12
$ws.Range("R141").Value = "Hortaliza"
